# Add a "Fall-2023" session label to each data row in column F
# (new column next to "Status"), and move the active selection to F8,
# matching the updated CHR_Excel_Sheet/Book2.xlsx layout.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2:F8").Value = "Fall-2023"

$ws.Range("F8").Select()
